$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, pushing the existing rows 19-23 down to 20-24.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly record.
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44511
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112022
$ws.Range("G19").Value = "Arveja Verde"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 15000
$ws.Range("L19").Value = 16000
$ws.Range("M19").Value = 15500
$ws.Range("N19").Value = "`$/saco 25 kilos"
$ws.Range("O19").Value = "Región del Maule"
$ws.Range("P19").Value = 620
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
